$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook

# Delete the worksheets that are no longer needed
$wb.Worksheets.Item("PurchaserDetails").Delete()
$wb.Worksheets.Item("Menu").Delete()
$wb.Worksheets.Item("Categories").Delete()

# Update the SignUp sheet: username/password value changes from chakk27 to chakk38
$ws = $wb.Worksheets.Item("SignUp")
$ws.Range("A2").Value = "chakk38"
$ws.Range("B2").Value = "chakk38"

# Make SignUp the active/selected sheet/tab
$ws.Activate()

$wb.Save()
